$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-28 04:59:08"
$wsZhCn.Range("G2").Value = "2016-01-28 04:59:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-28 04:59:19"
$wsDeDe.Range("G2").Value = "2016-01-28 05:00:06"
